$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal string value into a cell without Excel silently
# re-typing it as a number (which would strip things like trailing zeros,
# e.g. "108.00" -> 108). We flip the cell to Text, assign, then flip the
# format back to General so no stray number-format change is left behind.
function Set-TextValue($addr, $value) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.NumberFormat = "General"
}

# Apply the updated cell values row by row, matching the refreshed crypto feed.

# Row 2
$ws.Range("D2").Value = '33.951.43'
$ws.Range("E2").Value = '  -0.56%  '

# Row 3
$ws.Range("D3").Value = '1.778.63'
$ws.Range("E3").Value = '  -0.28%  '

# Row 4
$ws.Range("E4").Value = '  +0.19%  '

# Row 5
Set-TextValue "D5" '225.45'
$ws.Range("E5").Value = '  +1.81%  '

# Row 6
$ws.Range("E6").Value = '  +0.30%  '

# Row 7
$ws.Range("E7").Value = '  +0.15%  '

# Row 8
$ws.Range("E8").Value = '  +1.76%  '

# Row 9
Set-TextValue "D9" '0.289'
$ws.Range("E9").Value = '  +0.12%  '

# Row 10
Set-TextValue "D10" '0.0703'
$ws.Range("E10").Value = '  -0.89%  '

# Row 11
$ws.Range("E11").Value = '  +1.33%  '

# Row 12
$ws.Range("D12").Value = '2.036.22'
$ws.Range("E12").Value = '  -0.20%  '

# Row 13
$ws.Range("B13").Value = 'Chainlink'
$ws.Range("C13").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue "D13" '10.94'
$ws.Range("E13").Value = '  +2.55%  '

# Row 14
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.773.92'
$ws.Range("E14").Value = '  -0.74%  '

# Row 15
$ws.Range("E15").Value = '  -0.87%  '

# Row 16
$ws.Range("D16").Value = '33.930.89'
$ws.Range("E16").Value = '  -0.53%  '

# Row 17
$ws.Range("E17").Value = '  -2.04%  '

# Row 18
Set-TextValue "D18" '67.83'
$ws.Range("E18").Value = '  -0.18%  '

# Row 19
Set-TextValue "D19" '242.43'
$ws.Range("E19").Value = '  -1.03%  '

# Row 20
$ws.Range("E20").Value = '  +0.69%  '

# Row 21
$ws.Range("E21").Value = '  +0.17%  '

# Row 22
Set-TextValue "D22" '10.67'
$ws.Range("E22").Value = '  +0.15%  '

# Row 23
Set-TextValue "D23" '4.08'
$ws.Range("E23").Value = '  -0.61%  '

# Row 24
$ws.Range("E24").Value = '  -2.50%  '

# Row 25
Set-TextValue "D25" '160.19'
$ws.Range("E25").Value = '  +1.60%  '

# Row 26
Set-TextValue "D26" '16.26'
$ws.Range("E26").Value = '  -0.69%  '

# Row 27
Set-TextValue "D27" '7.07'
$ws.Range("E27").Value = '  +0.40%  '

# Row 28
$ws.Range("E28").Value = '  +0.33%  '

# Row 29
$ws.Range("E29").Value = '  +0.34%  '

# Row 30
$ws.Range("E30").Value = '  +3.00%  '

# Row 31
Set-TextValue "D31" '0.0511'
$ws.Range("E31").Value = '  -1.50%  '

# Row 32
$ws.Range("E32").Value = '  -1.44%  '

# Row 33
$ws.Range("E33").Value = '  -0.66%  '

# Row 34
$ws.Range("E34").Value = '  -1.64%  '

# Row 35
$ws.Range("D35").Value = '1.391.30'
$ws.Range("E35").Value = '  -0.21%  '

# Row 36
Set-TextValue "D36" '0.654'
$ws.Range("E36").Value = '  +4.02%  '

# Row 37
$ws.Range("E37").Value = '  -1.22%  '

# Row 38
$ws.Range("E38").Value = '  +0.21%  '

# Row 39
$ws.Range("E39").Value = '  +0.74%  '

# Row 40
$ws.Range("E40").Value = '  +4.48%  '

# Row 41
$ws.Range("E41").Value = '  -2.50%  '

# Row 42
$ws.Range("E42").Value = '  -3.64%  '

# Row 43
Set-TextValue "D43" '77.59'
$ws.Range("E43").Value = '  -2.44%  '

# Row 44
Set-TextValue "D44" '13.19'
$ws.Range("E44").Value = '  +11.73%  '

# Row 45
$ws.Range("E45").Value = '  +3.12%  '

# Row 46
$ws.Range("B46").Value = 'BabyDogeCoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue "D46" ([string]::Concat("0.0", [char]0x2086, "0136"))
$ws.Range("E46").Value = '  +15.86%  '

# Row 47
$ws.Range("B47").Value = 'Quant'
$ws.Range("C47").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue "D47" '108.00'
$ws.Range("E47").Value = '  +1.80%  '

# Row 48
$ws.Range("E48").Value = '  +0.23%  '

# Row 49
$ws.Range("E49").Value = '  -0.21%  '

# Row 50
$ws.Range("D50").Value = '1.935.15'
$ws.Range("E50").Value = '  +0.15%  '

# Row 51
$ws.Range("E51").Value = '  +0.61%  '
